$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the content of A4 (previously "SuperZitat 3")
$ws.Range("A4").ClearContents()

# Update the active selection to I9
$ws.Range("I9").Select()
